$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3000
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 3000
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H98").Value = 600.8889
$ws.Range("I98").Value = 600.8889
$ws.Range("K98").Value = 600.8889
$ws.Range("M98").Value = 897.1111
$ws.Range("H122").Value = 600.8889
$ws.Range("I122").Value = 600.8889
$ws.Range("K122").Value = 1802.6667
$ws.Range("M122").Value = 647.3332999999998
$ws.Range("H137").Value = 3347.6667
$ws.Range("I137").Value = 695.3333
$ws.Range("K137").Value = 2085.9999
$ws.Range("M137").Value = 464.0001000000002
$ws.Range("H138").Value = 4161.5137
$ws.Range("I138").Value = 2728.111
$ws.Range("J138").Value = 4622.25
$ws.Range("K138").Value = 8184.333
$ws.Range("L138").Value = 13866.75
$ws.Range("M138").Value = -3044.333
$ws.Range("N138").Value = -24146.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 1967.8
$ws.Range("I26").Value = 1967.8
$ws.Range("K26").Value = 1967.8
$ws.Range("M26").Value = -1637.8
$ws.Range("H32").Value = 6209.619
$ws.Range("I32").Value = 5021.2104
$ws.Range("K32").Value = 5021.2104
$ws.Range("M32").Value = -4734.2104
$ws.Range("H45").Value = 2209.1765
$ws.Range("J45").Value = 2571.2856
$ws.Range("L45").Value = 2571.2856
$ws.Range("N45").Value = -3325.2856
$ws.Range("H74").Value = 26798.8
$ws.Range("I74").Value = 39997.5
$ws.Range("J74").Value = 17999.666
$ws.Range("K74").Value = 39997.5
$ws.Range("L74").Value = 17999.666
$ws.Range("M74").Value = -39123.5
$ws.Range("N74").Value = -19747.666
$ws.Range("H77").Value = 26798.8
$ws.Range("I77").Value = 39997.5
$ws.Range("J77").Value = 17999.666
$ws.Range("K77").Value = 199987.5
$ws.Range("L77").Value = 89998.33
$ws.Range("M77").Value = -195619.5
$ws.Range("N77").Value = -98734.33
$ws.Range("H97").Value = 499.66666
$ws.Range("I97").Value = 499.66666
$ws.Range("K97").Value = 499.66666
$ws.Range("M97").Value = -3.666659999999979
$ws.Range("H132").Value = 4381.273
$ws.Range("I132").Value = 3033
$ws.Range("K132").Value = 9099
$ws.Range("M132").Value = -6569

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2681.8462
$ws.Range("I134").Value = 2681.8462
$ws.Range("K134").Value = 8045.5386
$ws.Range("M134").Value = -5510.5386

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H31").Value = 6744.3887
$ws.Range("I31").Value = 2692.75
$ws.Range("J31").Value = 14847.667
$ws.Range("K31").Value = 2692.75
$ws.Range("L31").Value = 14847.667
$ws.Range("M31").Value = -2397.75
$ws.Range("N31").Value = -15437.667
$ws.Range("H34").Value = 6744.3887
$ws.Range("I34").Value = 2692.75
$ws.Range("J34").Value = 14847.667
$ws.Range("K34").Value = 2692.75
$ws.Range("L34").Value = 14847.667
$ws.Range("M34").Value = -2490.75
$ws.Range("N34").Value = -15251.667
$ws.Range("H35").Value = 5000
$ws.Range("I35").Value = 5000
$ws.Range("K35").Value = 5000
$ws.Range("M35").Value = -4706
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 800.625
$ws.Range("J5").Value = 914
$ws.Range("L5").Value = 2742
$ws.Range("N5").Value = -2966
$ws.Range("H132").Value = 3846
$ws.Range("J132").Value = 5410
$ws.Range("L132").Value = 48690
$ws.Range("N132").Value = -53750
$ws.Range("H135").Value = 800.625
$ws.Range("J135").Value = 914
$ws.Range("L135").Value = 8226
$ws.Range("N135").Value = -13296
$ws.Range("H139").Value = 6507.5
$ws.Range("I139").Value = 5030
$ws.Range("K139").Value = 15090
$ws.Range("M139").Value = -9950

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3643.6
$ws.Range("I80").Value = 2498.4546
$ws.Range("K80").Value = 2498.4546
$ws.Range("M80").Value = -1500.4546
$ws.Range("H83").Value = 3643.6
$ws.Range("I83").Value = 2498.4546
$ws.Range("K83").Value = 12492.273
$ws.Range("M83").Value = -7500.273000000001
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 170
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 1099.2
$ws.Range("I122").Value = 999.25
$ws.Range("K122").Value = 2997.75
$ws.Range("M122").Value = -547.75
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 3717.5334
$ws.Range("I132").Value = 3160.6365
$ws.Range("K132").Value = 9481.9095
$ws.Range("M132").Value = -6951.9095

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3998.8
$ws.Range("J7").Value = 3998.3333
$ws.Range("L7").Value = 3998.3333
$ws.Range("N7").Value = -4222.3333
$ws.Range("H16").Value = 1104.2
$ws.Range("J16").Value = 497
$ws.Range("L16").Value = 497
$ws.Range("N16").Value = -837
$ws.Range("H22").Value = 3848.9312
$ws.Range("I22").Value = 3744.9375
$ws.Range("J22").Value = 3976.923
$ws.Range("K22").Value = 3744.9375
$ws.Range("L22").Value = 3976.923
$ws.Range("M22").Value = -3449.9375
$ws.Range("N22").Value = -4566.923
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H27").Value = 3848.9312
$ws.Range("I27").Value = 3744.9375
$ws.Range("J27").Value = 3976.923
$ws.Range("K27").Value = 3744.9375
$ws.Range("L27").Value = 3976.923
$ws.Range("M27").Value = -3637.9375
$ws.Range("N27").Value = -4190.923
$ws.Range("H61").Value = 4534.625
$ws.Range("I61").Value = 4468.4287
$ws.Range("J61").Value = 4998
$ws.Range("K61").Value = 4468.4287
$ws.Range("L61").Value = 4998
$ws.Range("M61").Value = -4266.4287
$ws.Range("N61").Value = -5402
$ws.Range("H113").Value = 4534.625
$ws.Range("I113").Value = 4468.4287
$ws.Range("J113").Value = 4998
$ws.Range("K113").Value = 4468.4287
$ws.Range("L113").Value = 4998
$ws.Range("M113").Value = -2298.4287
$ws.Range("N113").Value = -9338
$ws.Range("H126").Value = 3998.8
$ws.Range("J126").Value = 3998.3333
$ws.Range("L126").Value = 11994.9999
$ws.Range("N126").Value = -16934.9999
$ws.Range("H132").Value = 4810.1113
$ws.Range("I132").Value = 4383
$ws.Range("K132").Value = 13149
$ws.Range("M132").Value = -10619

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 956.8570999999999
$ws.Range("I81").Value = 956.8570999999999
$ws.Range("K81").Value = 1913.7142
$ws.Range("M81").Value = -852.7141999999999
$ws.Range("H84").Value = 956.8570999999999
$ws.Range("I84").Value = 956.8570999999999
$ws.Range("K84").Value = 9568.571
$ws.Range("M84").Value = -4264.571
$ws.Range("H113").Value = 737.4286
$ws.Range("I113").Value = 749
$ws.Range("K113").Value = 2247
$ws.Range("M113").Value = -77
$ws.Range("H132").Value = 2888.0938
$ws.Range("I132").Value = 2552.963
$ws.Range("K132").Value = 7658.889000000001
$ws.Range("M132").Value = -5128.889000000001
$ws.Range("H136").Value = 9266
$ws.Range("I136").Value = 9266
$ws.Range("K136").Value = 27798
$ws.Range("M136").Value = -25248
